$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value2 = 44505
$ws.Cells.Item(2, 10).Value2 = 440
$ws.Cells.Item(2, 11).Value2 = 900
$ws.Cells.Item(2, 12).Value2 = 1000
$ws.Cells.Item(2, 13).Value2 = 950
$ws.Cells.Item(2, 16).Value2 = 950

$ws.Cells.Item(3, 4).Value2 = 44545
$ws.Cells.Item(3, 10).Value2 = 4000
$ws.Cells.Item(3, 11).Value2 = 800
$ws.Cells.Item(3, 12).Value2 = 900
$ws.Cells.Item(3, 13).Value2 = 850
$ws.Cells.Item(3, 16).Value2 = 850

$ws.Cells.Item(4, 4).Value2 = 44511
$ws.Cells.Item(4, 10).Value2 = 500
$ws.Cells.Item(4, 11).Value2 = 900
$ws.Cells.Item(4, 12).Value2 = 1000
$ws.Cells.Item(4, 13).Value2 = 950
$ws.Cells.Item(4, 16).Value2 = 950

$ws.Cells.Item(5, 4).Value2 = 44512
$ws.Cells.Item(5, 10).Value2 = 600
$ws.Cells.Item(5, 11).Value2 = 900
$ws.Cells.Item(5, 12).Value2 = 1000
$ws.Cells.Item(5, 13).Value2 = 950
$ws.Cells.Item(5, 16).Value2 = 950

$ws.Cells.Item(6, 4).Value2 = 44518
$ws.Cells.Item(6, 10).Value2 = 400
$ws.Cells.Item(6, 11).Value2 = 800
$ws.Cells.Item(6, 12).Value2 = 900
$ws.Cells.Item(6, 13).Value2 = 850
$ws.Cells.Item(6, 16).Value2 = 850

$ws.Cells.Item(7, 4).Value2 = 44504
$ws.Cells.Item(7, 10).Value2 = 500
$ws.Cells.Item(7, 11).Value2 = 900
$ws.Cells.Item(7, 12).Value2 = 1000
$ws.Cells.Item(7, 13).Value2 = 950
$ws.Cells.Item(7, 16).Value2 = 950

$ws.Cells.Item(8, 4).Value2 = 44532
$ws.Cells.Item(8, 10).Value2 = 240
$ws.Cells.Item(8, 11).Value2 = 800
$ws.Cells.Item(8, 12).Value2 = 900
$ws.Cells.Item(8, 13).Value2 = 850
$ws.Cells.Item(8, 16).Value2 = 850

$ws.Cells.Item(9, 4).Value2 = 44510
$ws.Cells.Item(9, 10).Value2 = 600
$ws.Cells.Item(9, 11).Value2 = 900
$ws.Cells.Item(9, 12).Value2 = 1000
$ws.Cells.Item(9, 13).Value2 = 950
$ws.Cells.Item(9, 16).Value2 = 950

$ws.Cells.Item(10, 4).Value2 = 44524
$ws.Cells.Item(10, 10).Value2 = 400
$ws.Cells.Item(10, 11).Value2 = 800
$ws.Cells.Item(10, 12).Value2 = 900
$ws.Cells.Item(10, 13).Value2 = 850
$ws.Cells.Item(10, 16).Value2 = 850

$ws.Cells.Item(11, 4).Value2 = 44516
$ws.Cells.Item(11, 10).Value2 = 400
$ws.Cells.Item(11, 11).Value2 = 900
$ws.Cells.Item(11, 12).Value2 = 1000
$ws.Cells.Item(11, 13).Value2 = 950
$ws.Cells.Item(11, 16).Value2 = 950

$ws.Cells.Item(12, 4).Value2 = 44525
$ws.Cells.Item(12, 10).Value2 = 360
$ws.Cells.Item(12, 11).Value2 = 800
$ws.Cells.Item(12, 12).Value2 = 900
$ws.Cells.Item(12, 13).Value2 = 850
$ws.Cells.Item(12, 16).Value2 = 850

$ws.Cells.Item(13, 4).Value2 = 44553
$ws.Cells.Item(13, 10).Value2 = 8000
$ws.Cells.Item(13, 11).Value2 = 800
$ws.Cells.Item(13, 12).Value2 = 900
$ws.Cells.Item(13, 13).Value2 = 850
$ws.Cells.Item(13, 16).Value2 = 850

$ws.Cells.Item(14, 4).Value2 = 44503
$ws.Cells.Item(14, 10).Value2 = 400
$ws.Cells.Item(14, 11).Value2 = 900
$ws.Cells.Item(14, 12).Value2 = 1000
$ws.Cells.Item(14, 13).Value2 = 950
$ws.Cells.Item(14, 16).Value2 = 950

$ws.Cells.Item(15, 4).Value2 = 44530
$ws.Cells.Item(15, 10).Value2 = 300
$ws.Cells.Item(15, 11).Value2 = 800
$ws.Cells.Item(15, 12).Value2 = 900
$ws.Cells.Item(15, 13).Value2 = 850
$ws.Cells.Item(15, 16).Value2 = 850

$ws.Cells.Item(16, 4).Value2 = 44517
$ws.Cells.Item(16, 10).Value2 = 500
$ws.Cells.Item(16, 11).Value2 = 800
$ws.Cells.Item(16, 12).Value2 = 900
$ws.Cells.Item(16, 13).Value2 = 850
$ws.Cells.Item(16, 16).Value2 = 850

$ws.Cells.Item(17, 4).Value2 = 44476
$ws.Cells.Item(17, 10).Value2 = 300
$ws.Cells.Item(17, 11).Value2 = 1100
$ws.Cells.Item(17, 12).Value2 = 1200
$ws.Cells.Item(17, 13).Value2 = 1150
$ws.Cells.Item(17, 16).Value2 = 1150

$ws.Cells.Item(18, 4).Value2 = 44523
$ws.Cells.Item(18, 10).Value2 = 400
$ws.Cells.Item(18, 11).Value2 = 800
$ws.Cells.Item(18, 12).Value2 = 900
$ws.Cells.Item(18, 13).Value2 = 850
$ws.Cells.Item(18, 16).Value2 = 850

$ws.Cells.Item(19, 4).Value2 = 44537
$ws.Cells.Item(19, 10).Value2 = 400
$ws.Cells.Item(19, 11).Value2 = 800
$ws.Cells.Item(19, 12).Value2 = 900
$ws.Cells.Item(19, 13).Value2 = 850
$ws.Cells.Item(19, 16).Value2 = 850

$ws.Cells.Item(20, 4).Value2 = 44508
$ws.Cells.Item(20, 10).Value2 = 400
$ws.Cells.Item(20, 11).Value2 = 900
$ws.Cells.Item(20, 12).Value2 = 1000
$ws.Cells.Item(20, 13).Value2 = 950
$ws.Cells.Item(20, 16).Value2 = 950
